# 533-RBI-EPP-DB-SAR-REC-NON-RNI-CTPD-SAR-MD-TR-1-EarlyRePayment-Makerepayment1.xlsx
# "multi browser implementation for chrome and firefox and accounting cash"
#
# - Update a few interest/fee/due figures on the "Repayment Schedule" sheet
#   (rows 7-9) to reflect revised accounting-cash numbers.
# - Drop the stray empty P2 cell and replace it with an empty O2 cell on the
#   same sheet (same blank style).
# - Re-point the workbook's selections / active sheet: Input is no longer the
#   active tab, Summary's & Repayment Schedule's selections move, and
#   Transactions becomes the active tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Repayment Schedule: numeric corrections (rows 7, 8, 9)
# ---------------------------------------------------------------------------
$wsRepay = $wb.Worksheets.Item("Repayment Schedule")

$wsRepay.Range("H7").Value = 66.12
$wsRepay.Range("K7").Value = 899.45
$wsRepay.Range("P7").Value = 899.45

$wsRepay.Range("H8").Value = 57.79
$wsRepay.Range("K8").Value = 891.12
$wsRepay.Range("P8").Value = 891.12

$wsRepay.Range("H9").Value = 49.45
$wsRepay.Range("K9").Value = 882.78
$wsRepay.Range("P9").Value = 882.78

# Row 2: the empty "P2" cell moves to "O2" (same blank style).
$wsRepay.Range("P2").Clear()
$wsRepay.Range("O2").VerticalAlignment = -4108
$wsRepay.Range("O2").WrapText = $true

# ---------------------------------------------------------------------------
# Selections per-sheet (also reorders which tab ends up active/selected)
# ---------------------------------------------------------------------------

# Summary: selection moves from D18 to C4
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("C4").Select()

# Repayment Schedule: selection moves from A1:P1 to D8
$wsRepay.Range("D8").Select()

# Transactions becomes the active sheet/tab (was Input before).
$wsTransactions = $wb.Worksheets.Item("Transactions")
$wsTransactions.Activate()
